$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("excess_return_without_cost")
$ws2 = $wb.Worksheets.Item("excess_return_with_cost")

# Sheet: excess_return_without_cost - updated mean/std/annualized_return/information_ratio/max_drawdown values
$ws1.Range("C2").Value = 0.00076
$ws1.Range("C3").Value = 0.00137
$ws1.Range("D3").Value = 0.0012
$ws1.Range("C4").Value = 0.0009300000000000001
$ws1.Range("D4").Value = 0.00122
$ws1.Range("E4").Value = 0.00124
$ws1.Range("F4").Value = 0.00096
$ws1.Range("D5").Value = 0.0012
$ws1.Range("E5").Value = 0.00138
$ws1.Range("F5").Value = 0.00105
$ws1.Range("G5").Value = 0.00106
$ws1.Range("C6").Value = 0.00052
$ws1.Range("D6").Value = 0.00099
$ws1.Range("E6").Value = 0.00108
$ws1.Range("F6").Value = 0.0012
$ws1.Range("G6").Value = 0.00134
$ws1.Range("C7").Value = 0.00047
$ws1.Range("D7").Value = 0.00078
$ws1.Range("E7").Value = 0.00095
$ws1.Range("F7").Value = 0.00123
$ws1.Range("G7").Value = 0.00128
$ws1.Range("C8").Value = 0.02363
$ws1.Range("C9").Value = 0.01786
$ws1.Range("D9").Value = 0.01823
$ws1.Range("C10").Value = 0.01318
$ws1.Range("D10").Value = 0.01407
$ws1.Range("E10").Value = 0.01347
$ws1.Range("F10").Value = 0.01167
$ws1.Range("C11").Value = 0.01092
$ws1.Range("D11").Value = 0.01188
$ws1.Range("E11").Value = 0.01165
$ws1.Range("F11").Value = 0.01058
$ws1.Range("G11").Value = 0.01055
$ws1.Range("C12").Value = 0.01029
$ws1.Range("D12").Value = 0.0098
$ws1.Range("E12").Value = 0.01018
$ws1.Range("F12").Value = 0.01014
$ws1.Range("G12").Value = 0.01018
$ws1.Range("C13").Value = 0.009180000000000001
$ws1.Range("D13").Value = 0.008750000000000001
$ws1.Range("E13").Value = 0.008959999999999999
$ws1.Range("F13").Value = 0.00933
$ws1.Range("G13").Value = 0.00895
$ws1.Range("C14").Value = 0.182
$ws1.Range("C15").Value = 0.326
$ws1.Range("D15").Value = 0.28496
$ws1.Range("C16").Value = 0.22128
$ws1.Range("D16").Value = 0.29065
$ws1.Range("E16").Value = 0.29607
$ws1.Range("F16").Value = 0.2287
$ws1.Range("C17").Value = 0.13113
$ws1.Range("D17").Value = 0.28495
$ws1.Range("E17").Value = 0.32741
$ws1.Range("F17").Value = 0.25065
$ws1.Range("G17").Value = 0.25212
$ws1.Range("C18").Value = 0.12361
$ws1.Range("D18").Value = 0.23514
$ws1.Range("E18").Value = 0.25643
$ws1.Range("F18").Value = 0.28555
$ws1.Range("G18").Value = 0.31915
$ws1.Range("C19").Value = 0.11162
$ws1.Range("D19").Value = 0.18557
$ws1.Range("E19").Value = 0.22607
$ws1.Range("F19").Value = 0.29186
$ws1.Range("G19").Value = 0.30507
$ws1.Range("C20").Value = 0.49935
$ws1.Range("C21").Value = 1.1834
$ws1.Range("D21").Value = 1.01351
$ws1.Range("C22").Value = 1.08822
$ws1.Range("D22").Value = 1.33925
$ws1.Range("E22").Value = 1.42518
$ws1.Range("F22").Value = 1.27009
$ws1.Range("C23").Value = 0.77837
$ws1.Range("D23").Value = 1.55526
$ws1.Range("E23").Value = 1.82155
$ws1.Range("F23").Value = 1.5355
$ws1.Range("G23").Value = 1.54896
$ws1.Range("C24").Value = 0.7786
$ws1.Range("D24").Value = 1.55587
$ws1.Range("E24").Value = 1.63301
$ws1.Range("F24").Value = 1.82569
$ws1.Range("G24").Value = 2.03259
$ws1.Range("C25").Value = 0.78832
$ws1.Range("D25").Value = 1.37433
$ws1.Range("E25").Value = 1.63514
$ws1.Range("F25").Value = 2.02711
$ws1.Range("G25").Value = 2.20827
$ws1.Range("C26").Value = -0.53932
$ws1.Range("C27").Value = -0.27228
$ws1.Range("D27").Value = -0.34336
$ws1.Range("C28").Value = -0.2453
$ws1.Range("D28").Value = -0.27391
$ws1.Range("E28").Value = -0.2589
$ws1.Range("F28").Value = -0.29511
$ws1.Range("C29").Value = -0.19179
$ws1.Range("D29").Value = -0.18552
$ws1.Range("E29").Value = -0.16682
$ws1.Range("F29").Value = -0.19422
$ws1.Range("G29").Value = -0.17538
$ws1.Range("C30").Value = -0.21058
$ws1.Range("D30").Value = -0.18346
$ws1.Range("E30").Value = -0.1729
$ws1.Range("F30").Value = -0.116
$ws1.Range("G30").Value = -0.11739
$ws1.Range("C31").Value = -0.20559
$ws1.Range("D31").Value = -0.15441
$ws1.Range("E31").Value = -0.15447
$ws1.Range("F31").Value = -0.11283
$ws1.Range("G31").Value = -0.10983

# Sheet: excess_return_with_cost - updated mean/std/annualized_return/information_ratio/max_drawdown values
$ws2.Range("C2").Value = -0.00084
$ws2.Range("C3").Value = 0.00041
$ws2.Range("D3").Value = -0.00034
$ws2.Range("C4").Value = 0.00044
$ws2.Range("D4").Value = 0.00026
$ws2.Range("E4").Value = -0.00005
$ws2.Range("F4").Value = -0.00016
$ws2.Range("C5").Value = 0.00023
$ws2.Range("D5").Value = 0.00054
$ws2.Range("E5").Value = 0.00042
$ws2.Range("F5").Value = -0.00008000000000000001
$ws2.Range("G5").Value = -0.00015
$ws2.Range("C6").Value = 0.00026
$ws2.Range("D6").Value = 0.00046
$ws2.Range("E6").Value = 0.00032
$ws2.Range("F6").Value = 0.00023
$ws2.Range("G6").Value = 0.00021
$ws2.Range("C7").Value = 0.00027
$ws2.Range("D7").Value = 0.00036
$ws2.Range("E7").Value = 0.00034
$ws2.Range("F7").Value = 0.00043
$ws2.Range("G7").Value = 0.00036
$ws2.Range("C8").Value = 0.02365
$ws2.Range("C9").Value = 0.01786
$ws2.Range("D9").Value = 0.01824
$ws2.Range("C10").Value = 0.01318
$ws2.Range("D10").Value = 0.01407
$ws2.Range("E10").Value = 0.01347
$ws2.Range("F10").Value = 0.01167
$ws2.Range("C11").Value = 0.01092
$ws2.Range("D11").Value = 0.01188
$ws2.Range("E11").Value = 0.01166
$ws2.Range("F11").Value = 0.01059
$ws2.Range("G11").Value = 0.01056
$ws2.Range("C12").Value = 0.01029
$ws2.Range("D12").Value = 0.0098
$ws2.Range("E12").Value = 0.01018
$ws2.Range("F12").Value = 0.01014
$ws2.Range("G12").Value = 0.01018
$ws2.Range("C13").Value = 0.009180000000000001
$ws2.Range("D13").Value = 0.008750000000000001
$ws2.Range("E13").Value = 0.008959999999999999
$ws2.Range("F13").Value = 0.009339999999999999
$ws2.Range("G13").Value = 0.008959999999999999
$ws2.Range("C14").Value = -0.1997
$ws2.Range("C15").Value = 0.09846000000000001
$ws2.Range("D15").Value = -0.08185000000000001
$ws2.Range("C16").Value = 0.10506
$ws2.Range("D16").Value = 0.06136
$ws2.Range("E16").Value = -0.0118
$ws2.Range("F16").Value = -0.03864
$ws2.Range("C17").Value = 0.05374
$ws2.Range("D17").Value = 0.128
$ws2.Range("E17").Value = 0.09909
$ws2.Range("F17").Value = -0.01888
$ws2.Range("G17").Value = -0.03461
$ws2.Range("C18").Value = 0.06288000000000001
$ws2.Range("D18").Value = 0.11052
$ws2.Range("E18").Value = 0.07653
$ws2.Range("F18").Value = 0.05494
$ws2.Range("G18").Value = 0.04995
$ws2.Range("C19").Value = 0.06368
$ws2.Range("D19").Value = 0.08484
$ws2.Range("E19").Value = 0.08172
$ws2.Range("F19").Value = 0.1031
$ws2.Range("G19").Value = 0.08464000000000001
$ws2.Range("C20").Value = -0.5472399999999999
$ws2.Range("C21").Value = 0.35727
$ws2.Range("D21").Value = -0.29081
$ws2.Range("C22").Value = 0.5166500000000001
$ws2.Range("D22").Value = 0.28267
$ws2.Range("E22").Value = -0.0568
$ws2.Range("F22").Value = -0.21463
$ws2.Range("C23").Value = 0.319
$ws2.Range("D23").Value = 0.6986
$ws2.Range("E23").Value = 0.551
$ws2.Range("F23").Value = -0.11558
$ws2.Range("G23").Value = -0.21245
$ws2.Range("C24").Value = 0.39612
$ws2.Range("D24").Value = 0.73123
$ws2.Range("E24").Value = 0.4872
$ws2.Range("F24").Value = 0.35105
$ws2.Range("G24").Value = 0.31794
$ws2.Range("C25").Value = 0.44974
$ws2.Range("D25").Value = 0.62818
$ws2.Range("E25").Value = 0.59088
$ws2.Range("F25").Value = 0.71591
$ws2.Range("G25").Value = 0.61254
$ws2.Range("C26").Value = -1.45429
$ws2.Range("C27").Value = -0.35905
$ws2.Range("D27").Value = -0.91329
$ws2.Range("C28").Value = -0.29043
$ws2.Range("D28").Value = -0.46343
$ws2.Range("E28").Value = -0.64085
$ws2.Range("F28").Value = -0.63205
$ws2.Range("C29").Value = -0.23869
$ws2.Range("D29").Value = -0.29182
$ws2.Range("E29").Value = -0.35236
$ws2.Range("F29").Value = -0.47483
$ws2.Range("G29").Value = -0.49182
$ws2.Range("C30").Value = -0.23253
$ws2.Range("D30").Value = -0.23223
$ws2.Range("E30").Value = -0.30104
$ws2.Range("F30").Value = -0.2392
$ws2.Range("G30").Value = -0.26926
$ws2.Range("C31").Value = -0.22346
$ws2.Range("D31").Value = -0.19504
$ws2.Range("E31").Value = -0.2355
$ws2.Range("F31").Value = -0.18331
$ws2.Range("G31").Value = -0.18295
